# Fruta / hortaliza, semanal
# A new weekly price observation is inserted as row 104 (pushing the
# existing rows 104..206 down to 105..207), adding one more data row to
# the "Albahaca, Feria Lagunitas de Puerto Montt" series.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 104..206 down to 105..207, leaving a blank row 104 that
# inherits the formatting (incl. the date style) of the row it split from.
$ws.Rows(104).Insert()

# Populate the newly inserted row with the new observation.
$ws.Cells.Item(104, 1).Value = 4
$ws.Cells.Item(104, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(104, 3).Value = "Los Lagos"
$ws.Cells.Item(104, 4).Value = 45167
$ws.Cells.Item(104, 5).Value = 10
$ws.Cells.Item(104, 6).Value = 100112052
$ws.Cells.Item(104, 7).Value = "Albahaca"
$ws.Cells.Item(104, 8).Value = "Sin especificar"
$ws.Cells.Item(104, 9).Value = "Primera"
$ws.Cells.Item(104, 10).Value = 90
$ws.Cells.Item(104, 11).Value = 6000
$ws.Cells.Item(104, 12).Value = 6000
$ws.Cells.Item(104, 13).Value = 6000
$ws.Cells.Item(104, 14).Value = "`$/paquete"
$ws.Cells.Item(104, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(104, 16).Value = 6000
$ws.Cells.Item(104, 17).Value = 1
$ws.Cells.Item(104, 18).Value = "Hortaliza"
